$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.796.12'
$ws.Range('E2').Value = '  -0.24%  '
$ws.Range('D3').Value = '2.548.51'
$ws.Range('E3').Value = '  +0.51%  '
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '308.42'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -3.22%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '101.50'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +4.39%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.571'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.50%  '
$ws.Range('E8').Value = '  +0.06%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.534'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.73%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '36.34'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.42%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0809'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.01%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.40'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.95%  '
$ws.Range('E13').Value = '  -0.51%  '
$ws.Range('D14').Value = '2.943.98'
$ws.Range('E14').Value = '  +0.69%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.89'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +5.19%  '
$ws.Range('D16').Value = '2.545.12'
$ws.Range('E16').Value = '  -3.29%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.839'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.63%  '
$ws.Range('D18').Value = '42.840.55'
$ws.Range('E18').Value = '  -0.20%  '
$ws.Range('E19').Value = '  -1.93%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.39'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.43%  '
$ws.Range('D21').Value = '0.0₃0956'
$ws.Range('E21').Value = '  -1.23%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '69.23'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.71%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '247.12'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.43%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.91'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.59%  '
$ws.Range('E25').Value = '  +0.60%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '26.55'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.33%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '40.52'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -2.02%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.35'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -2.65%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '10.13'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -3.60%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '156.10'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.02%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.72'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -3.14%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0805'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.55%  '
$ws.Range('E34').Value = '  -3.31%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.28'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -2.16%  '
$ws.Range('E36').Value = '  -2.62%  '
$ws.Range('B37').Value = 'ApeXProtocol'
$ws.Range('C37').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.62'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +6.74%  '
$ws.Range('B38').Value = 'Celestia'
$ws.Range('C38').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '18.33'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -4.86%  '
$ws.Range('E39').Value = '  -0.95%  '
$ws.Range('E40').Value = '  -0.93%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.24'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +11.24%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '22.57'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +3.45%  '
$ws.Range('E43').Value = '  -0.14%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0301'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.67%  '
$ws.Range('E45').Value = '  +0.19%  '
$ws.Range('D46').Value = '1.985.27'
$ws.Range('E46').Value = '  -1.29%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '8.99'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.98%  '
$ws.Range('D48').Value = '2.796.85'
$ws.Range('E48').Value = '  +0.66%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '81.49'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -3.36%  '
$ws.Range('E50').Value = '  +0.66%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '73.67'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.28%  '
